$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 2")

# --- Correct / update the Database column text for a few rows ---
$ws.Range("D6").Value = "Database Design"
$ws.Range("D7").Value = "Database Design and implementation"
$ws.Range("D9").Value = "API `nBranching`nDatabase"

# --- Fill in "N/A" for Friday / obstacle columns (F,G) for rows 6-23 ---
for ($r = 6; $r -le 23; $r++) {
    $ws.Cells.Item($r, 6).Value = "N/A"
    $ws.Cells.Item($r, 7).Value = "N/A"
}

# --- Update view: scroll right one column and move selection to H23 ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H23").Select()
